# TC_169 - "Updated test data as per new implemenation"
# The Parent device for both IS device rows on the "Add IS Devices to EXI800"
# sheet is renamed from the generic "Exi800" to the specific "Exi800 - 1"
# instance, and the workbook is left with that sheet active/selected.

$wb = $excel.ActiveWorkbook

$wsAdd  = $wb.Worksheets.Item("Add EXI Devices Loop A")
$wsIs   = $wb.Worksheets.Item("Add IS Devices to EXI800")

# Parent column (J) for the two IS device rows now points at the specific
# Exi800 instance instead of the generic device name.
$wsIs.Range("J8").Value = "Exi800 - 1"
$wsIs.Range("J9").Value = "Exi800 - 1"

# Leave the workbook with "Add IS Devices to EXI800" as the active/selected
# sheet, with I9 selected.
$wsIs.Activate()
$wsIs.Range("I9").Select()
